# cryptos.xlsx refresh -- Sun Apr 16 06:25:11 UTC 2023 (GitHub Actions data pull)
# Refreshes the Price (D) and Volume(1h) (E) columns for each listed coin.
# Row 51 additionally drops out of the top-50 list and is replaced by a new coin
# (WOONetwork -> WEMIXTOKEN), so its Coin/Link/Price/Volume cells all change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Several Price values are numeric-looking text (e.g. 2.550, 0.6890) where a
# trailing zero is significant. Writing that straight to .Value lets Excel infer
# a Number type and the trailing zero is lost when it formats back (2.550 -> 2.55).
# Prefix Price values with an apostrophe -- same as typing '2.550 into the Excel UI --
# to force text storage and keep the exact digits, matching the source formatting.
function Set-PriceText($cellRef, $text) {
    $ws.Range($cellRef).Value = "`'$text"
}

Set-PriceText "D2" "30.592.28"
$ws.Range("E2").Value = "  +0.32%  "
Set-PriceText "D3" "2.116.74"
$ws.Range("E4").Value = "  +0.69%  "
Set-PriceText "D5" "336.95"
$ws.Range("E5").Value = "  +2.10%  "
Set-PriceText "D6" "1.009"
$ws.Range("E6").Value = "  +0.67%  "
Set-PriceText "D7" "0.5246"
$ws.Range("E7").Value = "  +0.67%  "
Set-PriceText "D8" "0.4559"
$ws.Range("E8").Value = "  +3.05%  "
Set-PriceText "D9" "54.75"
$ws.Range("E9").Value = "  +2.27%  "
Set-PriceText "D10" "0.09147"
$ws.Range("E10").Value = "  +2.42%  "
Set-PriceText "D11" "1.173"
$ws.Range("E11").Value = "  +1.81%  "
Set-PriceText "D12" "24.52"
$ws.Range("E12").Value = "  +1.09%  "
Set-PriceText "D13" "2.111.74"
$ws.Range("E13").Value = "  +0.66%  "
Set-PriceText "D14" "6.858"
$ws.Range("E14").Value = "  +2.52%  "
Set-PriceText "D15" "8.127"
$ws.Range("E15").Value = "  +5.79%  "
Set-PriceText "D16" "0.00001179"
$ws.Range("E16").Value = "  +5.01%  "
Set-PriceText "D17" "97.09"
Set-PriceText "D18" "1.011"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +1.27%  "
Set-PriceText "D20" "19.44"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  +0.70%  "
Set-PriceText "D22" "6.301"
$ws.Range("E22").Value = "  +0.79%  "
Set-PriceText "D23" "30.658.55"
$ws.Range("E23").Value = "  +0.42%  "
Set-PriceText "D24" "12.86"
$ws.Range("E24").Value = "  +4.46%  "
Set-PriceText "D25" "2.358"
$ws.Range("E25").Value = "  +1.89%  "
Set-PriceText "D26" "2.371.50"
$ws.Range("E26").Value = "  +1.21%  "
Set-PriceText "D27" "22.34"
$ws.Range("E27").Value = "  +0.38%  "
Set-PriceText "D28" "164.43"
$ws.Range("E28").Value = "  +0.40%  "
Set-PriceText "D29" "2.550"
$ws.Range("E29").Value = "  -0.23%  "
Set-PriceText "D30" "134.71"
$ws.Range("E30").Value = "  +2.44%  "
Set-PriceText "D31" "1.212"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("E32").Value = "  +0.56%  "
Set-PriceText "D33" "1.647"
$ws.Range("E33").Value = "  -0.25%  "
Set-PriceText "D34" "6.373"
$ws.Range("E34").Value = "  +3.54%  "
Set-PriceText "D35" "3.947"
$ws.Range("E35").Value = "  +1.10%  "
Set-PriceText "D36" "10.63"
$ws.Range("E36").Value = "  +5.90%  "
Set-PriceText "D37" "5.884"
$ws.Range("E37").Value = "  +7.45%  "
Set-PriceText "D38" "0.02636"
$ws.Range("E38").Value = "  +3.13%  "
Set-PriceText "D39" "0.06841"
$ws.Range("E39").Value = "  +0.39%  "
Set-PriceText "D40" "0.2329"
$ws.Range("E40").Value = "  +3.32%  "
$ws.Range("E41").Value = "  +0.10%  "
Set-PriceText "D42" "0.6890"
$ws.Range("E42").Value = "  +0.27%  "
Set-PriceText "D43" "1.260"
$ws.Range("E43").Value = "  +0.73%  "
Set-PriceText "D44" "14.71"
$ws.Range("E44").Value = "  +5.32%  "
Set-PriceText "D45" "0.6491"
$ws.Range("E45").Value = "  +2.72%  "
Set-PriceText "D46" "2.315"
$ws.Range("E46").Value = "  +5.40%  "
Set-PriceText "D47" "0.00000000367"
$ws.Range("E47").Value = "  +22.93%  "
Set-PriceText "D48" "3.691"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("E49").Value = "  +1.01%  "
Set-PriceText "D50" "83.42"
$ws.Range("E50").Value = "  +2.17%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-PriceText "D51" "1.183"
$ws.Range("E51").Value = "  -4.39%  "
